$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Date Analyse" column (H2 and H3) to the new analysis timestamp.
# Both rows now share the same date string, which collapses the previously
# distinct shared-string entries into a single one.
$ws.Range("H2").Value = "12/01/2026 09:56"
$ws.Range("H3").Value = "12/01/2026 09:56"
